# feat: add 2022-Q1 data
#
# The workbook currently has two sheets:
#   1) "2021-Q4" - per-fund holdings detail for 2021-Q4
#   2) "总计"     - quarter-by-quarter summary (date / count / market value)
#
# This script:
#   1) Renames the existing "总计" sheet to "2022-Q1" and rewrites it with the
#      2022-Q1 per-fund holdings detail (same shape as the "2021-Q4" sheet).
#   2) Inserts a brand new "总计" sheet right after "2022-Q1" that keeps the
#      summary table, with a new first data row for 2022-Q1 (and the old
#      2021-Q4 row pushed down).

$wb = $excel.ActiveWorkbook

$q4sheet = $wb.Worksheets.Item(1)      # "2021-Q4" (unchanged)
$q1sheet = $wb.Worksheets.Item(2)      # currently "总计" -> becomes "2022-Q1"

# --- 1) turn the old "总计" sheet into the new "2022-Q1" detail sheet -------

$q1sheet.Name = "2022-Q1"

# Header row is identical in shape/text/style to the "2021-Q4" sheet, so just
# copy it across (keeps the bold/centered/bordered header style intact).
$q4sheet.Range("B1:H1").Copy($q1sheet.Range("B1:H1"))

# Row 2 - 501030
$c = $q1sheet.Cells.Item(2, 1); $c.Value = 0
$c = $q1sheet.Cells.Item(2, 2); $c.NumberFormat = "@"; $c.Value = "501030"
$c = $q1sheet.Cells.Item(2, 3); $c.NumberFormat = "@"; $c.Value = "汇添富中证环境治理指数（LOF）A"
$c = $q1sheet.Cells.Item(2, 4); $c.NumberFormat = "@"; $c.Value = "6.61"
$c = $q1sheet.Cells.Item(2, 5); $c.NumberFormat = "@"; $c.Value = "93.20"
$c = $q1sheet.Cells.Item(2, 6); $c.NumberFormat = "@"; $c.Value = "2.38"
$c = $q1sheet.Cells.Item(2, 7); $c.NumberFormat = "@"; $c.Value = "0.1573"
$c = $q1sheet.Cells.Item(2, 8); $c.Value = 3

# Row 3 - 501031
$c = $q1sheet.Cells.Item(3, 1); $c.Value = 1
$c = $q1sheet.Cells.Item(3, 2); $c.NumberFormat = "@"; $c.Value = "501031"
$c = $q1sheet.Cells.Item(3, 3); $c.NumberFormat = "@"; $c.Value = "汇添富中证环境治理指数（LOF）C"
$c = $q1sheet.Cells.Item(3, 4); $c.NumberFormat = "@"; $c.Value = "2.74"
$c = $q1sheet.Cells.Item(3, 5); $c.NumberFormat = "@"; $c.Value = "93.20"
$c = $q1sheet.Cells.Item(3, 6); $c.NumberFormat = "@"; $c.Value = "2.38"
$c = $q1sheet.Cells.Item(3, 7); $c.NumberFormat = "@"; $c.Value = "0.0652"
$c = $q1sheet.Cells.Item(3, 8); $c.Value = 3

# Row 4 - 164908
$c = $q1sheet.Cells.Item(4, 1); $c.Value = 2
$c = $q1sheet.Cells.Item(4, 2); $c.NumberFormat = "@"; $c.Value = "164908"
$c = $q1sheet.Cells.Item(4, 3); $c.NumberFormat = "@"; $c.Value = "交银施罗德中证环境治理指数（LOF）"
$c = $q1sheet.Cells.Item(4, 4); $c.NumberFormat = "@"; $c.Value = "2.12"
$c = $q1sheet.Cells.Item(4, 5); $c.NumberFormat = "@"; $c.Value = "93.72"
$c = $q1sheet.Cells.Item(4, 6); $c.NumberFormat = "@"; $c.Value = "2.43"
$c = $q1sheet.Cells.Item(4, 7); $c.NumberFormat = "@"; $c.Value = "0.0515"
$c = $q1sheet.Cells.Item(4, 8); $c.Value = 1

# Column A (index) cells use the same bold/centered/bordered style as the
# header column in the "2021-Q4" sheet.
$q4sheet.Range("A2:A4").Copy($q1sheet.Range("A2:A4"))
$q1sheet.Cells.Item(2, 1).Value = 0
$q1sheet.Cells.Item(3, 1).Value = 1
$q1sheet.Cells.Item(4, 1).Value = 2

# --- 2) insert a fresh "总计" sheet right after "2022-Q1" -------------------

$totalSheet = $wb.Worksheets.Add($null, $q1sheet)
$totalSheet.Name = "总计"

# Header row: 日期 / 持有数量(只) / 持有市值(亿元)
$c = $totalSheet.Cells.Item(1, 2); $c.NumberFormat = "@"; $c.Value = "日期"
$c = $totalSheet.Cells.Item(1, 3); $c.NumberFormat = "@"; $c.Value = "持有数量(只)"
$c = $totalSheet.Cells.Item(1, 4); $c.NumberFormat = "@"; $c.Value = "持有市值(亿元)"

# New row for 2022-Q1
$totalSheet.Cells.Item(2, 1).Value = 0
$c = $totalSheet.Cells.Item(2, 2); $c.NumberFormat = "@"; $c.Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 3
$totalSheet.Cells.Item(2, 4).Value = 0.27

# Old row, now pushed to row 3
$totalSheet.Cells.Item(3, 1).Value = 1
$c = $totalSheet.Cells.Item(3, 2); $c.NumberFormat = "@"; $c.Value = "2021-Q4"
$totalSheet.Cells.Item(3, 3).Value = 3
$totalSheet.Cells.Item(3, 4).Value = 0.23

# Apply the same header style (bold/centered/bordered) used across the rest
# of the workbook to the new header row + the "A" index column.
$q4sheet.Range("B1:D1").Copy()
$totalSheet.Range("B1:D1").PasteSpecial(-4122) # xlPasteFormats
$q4sheet.Range("A2:A2").Copy()
$totalSheet.Range("A2:A3").PasteSpecial(-4122) # xlPasteFormats

$totalSheet.Cells.Item(1,1).Select()
